$wb = $excel.ActiveWorkbook

# Sheet "展览" (rId1 / sheet1.xml) - update "想去人数" (column F) counts
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 75
$wsExhibit.Range("F4").Value = 235
$wsExhibit.Range("F6").Value = 9948
$wsExhibit.Range("F9").Value = 1241
$wsExhibit.Range("F10").Value = 4982
$wsExhibit.Range("F14").Value = 118
$wsExhibit.Range("F17").Value = 291
$wsExhibit.Range("F19").Value = 112
$wsExhibit.Range("F21").Value = 1492

# Sheet "全部类型" (rId4 / sheet4.xml) - same events, same updated counts
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 75
$wsAll.Range("F5").Value = 235
$wsAll.Range("F7").Value = 9948
$wsAll.Range("F10").Value = 1241
$wsAll.Range("F11").Value = 4982
$wsAll.Range("F15").Value = 118
$wsAll.Range("F18").Value = 291
$wsAll.Range("F20").Value = 112
$wsAll.Range("F22").Value = 1492
